$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the data rows 2-10 (columns B:G) down into rows 3-11, dropping the
# old row 11 data, then populate row 2 with the new observation's values.
# Process from the bottom up so each source row is read before being
# overwritten by the row above it.

$ws.Range("B11:G11").Value2 = $ws.Range("B10:G10").Value2
$ws.Range("B10:G10").Value2 = $ws.Range("B9:G9").Value2
$ws.Range("B9:G9").Value2   = $ws.Range("B8:G8").Value2
$ws.Range("B8:G8").Value2   = $ws.Range("B7:G7").Value2
$ws.Range("B7:G7").Value2   = $ws.Range("B6:G6").Value2
$ws.Range("B6:G6").Value2   = $ws.Range("B5:G5").Value2
$ws.Range("B5:G5").Value2   = $ws.Range("B4:G4").Value2
$ws.Range("B4:G4").Value2   = $ws.Range("B3:G3").Value2
$ws.Range("B3:G3").Value2   = $ws.Range("B2:G2").Value2

$ws.Range("B2").Value2 = 0.08266386729847572
$ws.Range("C2").Value2 = 1.266710845429791
$ws.Range("D2").Value2 = 10.58350078540567
$ws.Range("E2").Value2 = 3.253229285710688
$ws.Range("F2").Value2 = 3.288115649630924
$ws.Range("G2").Value2 = 46
